$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '28.477.43'
$ws.Range("E2").Value = '  +0.89%  '

$ws.Range("D3").Value = '1.871.43'
$ws.Range("E3").Value = '  +0.76%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("E4").Value = '  -0.57%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.85'
$ws.Range("E5").Value = '  +0.53%  '

$ws.Range("E6").Value = '  -0.43%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5092'
$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3897'
$ws.Range("E8").Value = '  -0.44%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08329'
$ws.Range("E9").Value = '  +0.76%  '

$ws.Range("E10").Value = '  -0.58%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.73'
$ws.Range("E11").Value = '  +0.07%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.227'
$ws.Range("E12").Value = '  +0.43%  '

$ws.Range("D13").Value = '1.870.14'
$ws.Range("E13").Value = '  +0.99%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.36'
$ws.Range("E14").Value = '  +0.60%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.267'
$ws.Range("E15").Value = '  +1.06%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.010'
$ws.Range("E16").Value = '  -0.40%  '

$ws.Range("E17").Value = '  +0.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.13'
$ws.Range("E18").Value = '  +0.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06729'
$ws.Range("E19").Value = '  +0.66%  '

$ws.Range("E20").Value = '  +0.97%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.006'
$ws.Range("E21").Value = '  -0.49%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.904'
$ws.Range("E22").Value = '  -0.22%  '

$ws.Range("D23").Value = '28.487.16'
$ws.Range("E23").Value = '  +0.91%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.11'
$ws.Range("E24").Value = '  +0.50%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.225'
$ws.Range("E25").Value = '  -1.01%  '

$ws.Range("D26").Value = '2.086.24'
$ws.Range("E26").Value = '  +1.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.90'
$ws.Range("E27").Value = '  +1.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.61'
$ws.Range("E28").Value = '  +0.08%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.401'
$ws.Range("E29").Value = '  +1.74%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.16'
$ws.Range("E30").Value = '  +0.00%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1042'
$ws.Range("E31").Value = '  -0.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.035'
$ws.Range("E32").Value = '  +1.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.744'
$ws.Range("E33").Value = '  -0.64%  '

$ws.Range("E34").Value = '  -0.66%  '

$ws.Range("E35").Value = '  +1.34%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06545'
$ws.Range("E36").Value = '  +1.78%  '

$ws.Range("E37").Value = '  -0.44%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.822'
$ws.Range("E38").Value = '  -2.55%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.021'
$ws.Range("E39").Value = '  +2.11%  '

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.245'
$ws.Range("E40").Value = '  -0.02%  '

$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.184'
$ws.Range("E41").Value = '  +0.61%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6367'
$ws.Range("E42").Value = '  -0.33%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.06'
$ws.Range("E43").Value = '  +0.10%  '

$ws.Range("E44").Value = '  -0.38%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5988'
$ws.Range("E45").Value = '  +0.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.02'
$ws.Range("E46").Value = '  +2.14%  '

$ws.Range("E47").Value = '  -0.11%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.001'
$ws.Range("E48").Value = '  +1.59%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.215'
$ws.Range("E49").Value = '  +1.16%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '121.76'
$ws.Range("E50").Value = '  +0.83%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.149'
$ws.Range("E51").Value = '  -10.07%  '
